$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. "549.22", "1.00", "0.0000267") are stored as text, matching
# the source data which uses inline strings throughout.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.572.97"
$ws.Range("E2").Value = "  -3.58%  "

$ws.Range("D3").Value = "3.315.08"
$ws.Range("E3").Value = "  -5.10%  "

$ws.Range("E4").Value = "  +0.41%  "

$ws.Range("D5").Value = "549.22"
$ws.Range("E5").Value = "  -4.21%  "

$ws.Range("D6").Value = "172.84"
$ws.Range("E6").Value = "  -2.69%  "

$ws.Range("D7").Value = "0.610"
$ws.Range("E7").Value = "  -3.66%  "

$ws.Range("D9").Value = "3.306.30"
$ws.Range("E9").Value = "  -5.30%  "

$ws.Range("D10").Value = "0.616"
$ws.Range("E10").Value = "  -2.17%  "

$ws.Range("D11").Value = "0.157"
$ws.Range("E11").Value = "  +1.07%  "

$ws.Range("D12").Value = "53.46"
$ws.Range("E12").Value = "  -1.26%  "

$ws.Range("D13").Value = "0.0000267"
$ws.Range("E13").Value = "  -1.32%  "

$ws.Range("D14").Value = "8.93"
$ws.Range("E14").Value = "  -2.53%  "

$ws.Range("D15").Value = "3.860.04"
$ws.Range("E15").Value = "  -4.67%  "

$ws.Range("D16").Value = "18.03"
$ws.Range("E16").Value = "  -1.05%  "

$ws.Range("D17").Value = "0.117"
$ws.Range("E17").Value = "  -3.35%  "

$ws.Range("D18").Value = "3.329.42"
$ws.Range("E18").Value = "  -4.61%  "

$ws.Range("D19").Value = "11.70"
$ws.Range("E19").Value = "  -2.93%  "

$ws.Range("D20").Value = "63.579.98"
$ws.Range("E20").Value = "  -3.33%  "

$ws.Range("D21").Value = "0.964"
$ws.Range("E21").Value = "  -3.48%  "

$ws.Range("D22").Value = "423.08"
$ws.Range("E22").Value = "  +2.53%  "

$ws.Range("D23").Value = "4.62"
$ws.Range("E23").Value = "  +9.27%  "

$ws.Range("D24").Value = "4.07"
$ws.Range("E24").Value = "  -2.66%  "

$ws.Range("D25").Value = "83.48"
$ws.Range("E25").Value = "  -1.80%  "

$ws.Range("D26").Value = "12.89"
$ws.Range("E26").Value = "  +1.78%  "

$ws.Range("D27").Value = "10.51"
$ws.Range("E27").Value = "  -3.40%  "

$ws.Range("D28").Value = "2.79"
$ws.Range("E28").Value = "  -1.41%  "

$ws.Range("D29").Value = "8.61"
$ws.Range("E29").Value = "  -3.91%  "

$ws.Range("D30").Value = "29.32"
$ws.Range("E30").Value = "  -2.81%  "

$ws.Range("D31").Value = "6.49"
$ws.Range("E31").Value = "  +2.27%  "

$ws.Range("D32").Value = "585.44"
$ws.Range("E32").Value = "  -5.21%  "

$ws.Range("D33").Value = "11.33"
$ws.Range("E33").Value = "  -2.32%  "

$ws.Range("D34").Value = "0.106"
$ws.Range("E34").Value = "  -2.78%  "

$ws.Range("D35").Value = "58.34"
$ws.Range("E35").Value = "  -1.87%  "

$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").Value = "0.142"
$ws.Range("E37").Value = "  -5.91%  "

$ws.Range("D38").Value = "3.43"
$ws.Range("E38").Value = "  +2.66%  "

$ws.Range("D39").Value = "35.09"
$ws.Range("E39").Value = "  -5.29%  "

$ws.Range("D40").Value = "0.0₃0739"
$ws.Range("E40").Value = "  -7.07%  "

$ws.Range("D41").Value = "0.361"
$ws.Range("E41").Value = "  -4.45%  "

$ws.Range("D42").Value = "3.083.90"
$ws.Range("E42").Value = "  -6.82%  "

$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.60%  "

$ws.Range("D44").Value = "2.80"
$ws.Range("E44").Value = "  -2.53%  "

$ws.Range("D45").Value = "3.20"
$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("D46").Value = "0.0402"
$ws.Range("E46").Value = "  -2.59%  "

$ws.Range("D47").Value = "2.41"
$ws.Range("E47").Value = "  -3.66%  "

$ws.Range("D48").Value = "0.128"
$ws.Range("E48").Value = "  -2.83%  "

$ws.Range("E49").Value = "  -4.07%  "

$ws.Range("D50").Value = "132.46"
$ws.Range("E50").Value = "  -3.98%  "

$ws.Range("D51").Value = "8.09"
$ws.Range("E51").Value = "  -3.45%  "
